$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing item row (row 7), pushing the
# PANTOLOC item (and the totals/footer rows below it) down by one row.
$ws.Rows.Item(7).Insert()

# The inserted row has no formatting of its own; clone the exact
# look (styles + borders + fills) of the item row directly below it
# (which now holds the PANTOLOC data that used to live in row 7) so the
# new row matches the report's item-row style precisely.
$ws.Range("A8:Q8").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)

# Recreate the same cell merges used by every item row.
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()

# Populate the new item (#1: LEZBERG TRIO) in row 7.
$ws.Range("A7").Value2 = 1
$ws.Range("C7").Value2 = "LEZBERG TRIO 20/5/12.5 TAB"
$ws.Range("H7").Value2 = "0:2"

# L7 and P7 hold numeric-looking text ("0" / "37.6200") that must stay
# text cells (matching how the source report stores every item field as
# text) rather than being auto-coerced to numbers. Flip to a text format
# while assigning, then restore the original numeric display format so
# the cell style id is unchanged.
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value2 = "0"
$ws.Range("L7").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N7").Value2 = "114.00"

$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value2 = "37.6200"
$ws.Range("P7").NumberFormat = "0.00"

$ws.Range("Q7").Value2 = "0:1"

# Row 7 (the new LEZBERG row) gets its own autofit-derived height.
$ws.Rows.Item(7).RowHeight = 24.75

# Renumber the item that shifted from row 7 to row 8 (PANTOLOC is now
# item #2).
$ws.Range("A8").Value2 = 2

# The totals row (now row 9) picks up the new item's price, and its
# rendered height changes slightly with the new total.
$ws.Range("N9").Value2 = 88.62
$ws.Rows.Item(9).RowHeight = 26.25
